$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Update the "datetimeFigureOut" Date placeholder text on the Slide Master
#    and on every Slide Layout: 10/26/2017 -> 10/27/2017.
# ---------------------------------------------------------------------------
$master = $p.Slides.Item(1).Master
$masterDate = $master.Shapes.Item(4)
$masterDate.TextFrame.TextRange.Text = "10/27/2017"

for ($l = 1; $l -le $master.CustomLayouts.Count; $l++) {
    $layout = $master.CustomLayouts.Item($l)
    $layoutDate = $layout.Shapes.Item(4)
    $layoutDate.TextFrame.TextRange.Text = "10/27/2017"
}
